$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Drop the redundant "Unnamed: 0" index column (old column B). Everything
#    to its right (Country .. hypothetical predictions) shifts one column
#    left automatically.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).Delete()

# ---------------------------------------------------------------------------
# 2. Drop the old "Neural network Predictions hypothetical test countries"
#    column. After step 1 it now sits at column J (10th column); removing it
#    shifts nothing further right, it's the last used column.
# ---------------------------------------------------------------------------
$ws.Columns.Item(10).Delete()

# ---------------------------------------------------------------------------
# 3. Header for column I now must read the old "test countries" title
#    (it previously held "Pop*1.1" before the shift -- after the column B
#    delete it already correctly reads "Neural network Predictions test
#    countries" coming from the old J1, so nothing else to rename here).
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 4. Gini Index values (column C) switch from whole percentages (32.8) to
#    fractional values (0.328).
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = 0.328
$ws.Range("C3").Value = 0.357
$ws.Range("C4").Value = 0.327
$ws.Range("C5").Value = 0.408
$ws.Range("C6").Value = 0.435
$ws.Range("C7").Value = 0.571
$ws.Range("C8").Value = 0.42
$ws.Range("C9").Value = 0.39
$ws.Range("C10").Value = 0.457
$ws.Range("C11").Value = 0.457

# ---------------------------------------------------------------------------
# 5. Recomputed "Neural network Predictions test countries" values
#    (column I).
# ---------------------------------------------------------------------------
$ws.Range("I2").Value = 422547.5602236092
$ws.Range("I3").Value = 1517195.957459688
$ws.Range("I4").Value = 112809.336062327
$ws.Range("I5").Value = 1130416.920038164
$ws.Range("I6").Value = 1083931.953636229
$ws.Range("I7").Value = 1173707.549841166
$ws.Range("I8").Value = 1704420.667387366
$ws.Range("I9").Value = 527852.0890156627
$ws.Range("I10").Value = 731229.9269056916
$ws.Range("I11").Value = 929388.7327637076

# ---------------------------------------------------------------------------
# 6. New "Neural network Predictions hypothetical test countries" column
#    (column J), header + recomputed values. Column J no longer exists after
#    the deletes above, so first clone the header formatting from I1 (bold,
#    centered, bordered) before writing the new title text.
# ---------------------------------------------------------------------------
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "Neural network Predictions hypothetical test countries"
$ws.Range("J2").Value = 697850.3890247345
$ws.Range("J3").Value = 2387619.971201897
$ws.Range("J4").Value = 19335.79575616494
$ws.Range("J5").Value = 1303562.652504861
$ws.Range("J6").Value = 823892.7610159516
$ws.Range("J7").Value = 458444.4515973032
$ws.Range("J8").Value = 2204333.234232783
$ws.Range("J9").Value = 137071.0585970879
$ws.Range("J10").Value = 175536.5894390792
$ws.Range("J11").Value = 465069.2912610471

# ---------------------------------------------------------------------------
# 7. New row for Egypt, the 11th test country (index 10).
# ---------------------------------------------------------------------------
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Egypt"
$ws.Range("C12").Value = 0.315
$ws.Range("D12").Value = 100400000
$ws.Range("E12").Value = 43
$ws.Range("F12").Value = 2824316
$ws.Range("G12").Value = 258407
$ws.Range("H12").Value = 110440000
$ws.Range("I12").Value = 1437885.004332185
$ws.Range("J12").Value = 2497490.042134404

$excel.CutCopyMode = 0
